$wb = $excel.ActiveWorkbook

# --- Sheet "Todos": add DUE_DATE (G2) and CREATED_DATE (H2) ---
$todos = $wb.Worksheets.Item("Todos")

# Set H2 first and give it a real date number format (builtin mm-dd-yy / numFmtId 14).
$todos.Range("H2").NumberFormatLocal = "mm-dd-yy"
$todos.Range("H2").Value = 45805

# Copy H2's format (only) onto G2 so both cells share a single style entry,
# then set G2's value.
$todos.Range("H2").Copy()
$todos.Range("G2").PasteSpecial(-4122)
$todos.Range("G2").Value = 45838

[void]$todos.Range("F7").Select()

# --- Sheet "Updates": add UPDATE_DATE (B2) ---
$updates = $wb.Worksheets.Item("Updates")

$todos.Range("H2").Copy()
$updates.Range("B2").PasteSpecial(-4122)
$updates.Range("B2").Value = 45835

[void]$updates.Range("B5").Select()
